$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 4 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Wil je dit oppakken?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #2: Wil je dit oppakken?"
$logs.Range("D4").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E4").Value = "Beste klant,`nDank u voor uw e-mail. Om u zo goed mogelijk van dienst te kunnen zijn, wil ik u vragen om meer informatie te verstrekken over de specifieke kwestie die u wilt dat we oppakken. Kunt u mij wat meer details geven over wat u precies nodig heeft? Op die manier kunnen we u beter helpen.`nMet vriendelijke groet,`n[Naam van de e-mailassistent]`n[Naam van het bedrijf]"
$logs.Range("F4").Value = "2025-07-27 19:14:34"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Nee"
$logs.Range("I4").Value = "Ja"
$logs.Range("J4").Value = "Nee"
$logs.Rows.Item(4).AutoFit()

# --- Dashboard sheet: append row 3 ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B3").Value = 1

# --- Update chart series references to include the new Dashboard row ---
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$3"
$series.Values = "='Dashboard'!`$B`$2:`$B`$3"

# --- Expand conditional formatting ranges on Logs sheet to include row 4 ---
$ranges = @("D2:D3", "G2:G3", "H2:H3", "I2:I3", "J2:J3")
$newRanges = @("D2:D4", "G2:G4", "H2:H4", "I2:I4", "J2:J4")
for ($i = 0; $i -lt $ranges.Length; $i++) {
    $fcs = $logs.Range($ranges[$i]).FormatConditions
    for ($j = 1; $j -le $fcs.Count; $j++) {
        $fcs.Item($j).ModifyAppliesToRange($logs.Range($newRanges[$i]))
    }
}
